# Applies the 2026-02-02 Betfair Back/Lay odds-sheet update to Sheet1:
#  - Refreshed odds across many existing match rows (2-30)
#  - Row 17 / Row 18: Home & Away team names were swapped between the two
#    Italian Serie C fixtures (Giana Erminio vs Union Brescia <-> Renate vs
#    ASD Alcione), with each row's odds refreshed to match the new fixture
#  - Three brand new match rows appended at the bottom: 31 & 32 (Argentinian
#    Primera Division) and 33 (Colombian Primera A) -- dimension grows to AO33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Scattered odds/name updates on existing rows (2-30) ----

# Row 2
$ws.Cells.Item(2, 6).Value = 10
$ws.Cells.Item(2, 8).Value = 1.22
$ws.Cells.Item(2, 9).Value = 1.27
$ws.Cells.Item(2, 17).Value = 1.35
$ws.Cells.Item(2, 19).Value = 1.88
$ws.Cells.Item(2, 20).Value = 1.86
$ws.Cells.Item(2, 21).Value = 1.94
$ws.Cells.Item(2, 22).Value = 4.6
$ws.Cells.Item(2, 24).Value = 46
$ws.Cells.Item(2, 25).Value = 970
$ws.Cells.Item(2, 31).Value = 970
$ws.Cells.Item(2, 35).Value = 38

# Row 3
$ws.Cells.Item(3, 7).Value = 2.94
$ws.Cells.Item(3, 9).Value = 3.4
$ws.Cells.Item(3, 10).Value = 3.1
$ws.Cells.Item(3, 11).Value = 3.5
$ws.Cells.Item(3, 12).Value = 1.41
$ws.Cells.Item(3, 13).Value = 1.08
$ws.Cells.Item(3, 14).Value = 2.94
$ws.Cells.Item(3, 18).Value = 1.24
$ws.Cells.Item(3, 19).Value = 4.3
$ws.Cells.Item(3, 20).Value = 1.89
$ws.Cells.Item(3, 21).Value = 1.89
$ws.Cells.Item(3, 23).Value = 1.52
$ws.Cells.Item(3, 24).Value = 13
$ws.Cells.Item(3, 25).Value = 12.5
$ws.Cells.Item(3, 26).Value = 25
$ws.Cells.Item(3, 27).Value = 70
$ws.Cells.Item(3, 28).Value = 11
$ws.Cells.Item(3, 29).Value = 8.800000000000001
$ws.Cells.Item(3, 30).Value = 17
$ws.Cells.Item(3, 31).Value = 55
$ws.Cells.Item(3, 32).Value = 21
$ws.Cells.Item(3, 33).Value = 15
$ws.Cells.Item(3, 34).Value = 24
$ws.Cells.Item(3, 35).Value = 75
$ws.Cells.Item(3, 36).Value = 55
$ws.Cells.Item(3, 37).Value = 42
$ws.Cells.Item(3, 38).Value = 65
$ws.Cells.Item(3, 39).Value = 160
$ws.Cells.Item(3, 40).Value = 42
$ws.Cells.Item(3, 41).Value = 60

# Row 4
$ws.Cells.Item(4, 7).Value = 1.44
$ws.Cells.Item(4, 23).Value = 3.25

# Row 6
$ws.Cells.Item(6, 20).Value = 1.52
$ws.Cells.Item(6, 21).Value = 2.62

# Row 7
$ws.Cells.Item(7, 17).Value = 1.89
$ws.Cells.Item(7, 19).Value = 1.9

# Row 8
$ws.Cells.Item(8, 6).Value = 4.1
$ws.Cells.Item(8, 7).Value = 6.6
$ws.Cells.Item(8, 8).Value = 1.62
$ws.Cells.Item(8, 9).Value = 1.78
$ws.Cells.Item(8, 11).Value = 5.8
$ws.Cells.Item(8, 14).Value = 2.4
$ws.Cells.Item(8, 16).Value = 2.4
$ws.Cells.Item(8, 17).Value = 1.5
$ws.Cells.Item(8, 22).Value = 2.28
$ws.Cells.Item(8, 23).Value = 1.17

# Row 12
$ws.Cells.Item(12, 6).Value = 1.99
$ws.Cells.Item(12, 7).Value = 2.24
$ws.Cells.Item(12, 8).Value = 3.6
$ws.Cells.Item(12, 9).Value = 4.5
$ws.Cells.Item(12, 10).Value = 3.4
$ws.Cells.Item(12, 11).Value = 4.2
$ws.Cells.Item(12, 22).Value = 1.29
$ws.Cells.Item(12, 23).Value = 1.8

# Row 13
$ws.Cells.Item(13, 6).Value = 3.9
$ws.Cells.Item(13, 7).Value = 4.8
$ws.Cells.Item(13, 8).Value = 1.8
$ws.Cells.Item(13, 11).Value = 4.7
$ws.Cells.Item(13, 23).Value = 1.26
$ws.Cells.Item(13, 32).Value = 42
$ws.Cells.Item(13, 36).Value = 80
$ws.Cells.Item(13, 38).Value = 44
$ws.Cells.Item(13, 40).Value = 28

# Row 15
$ws.Cells.Item(15, 18).Value = 2
$ws.Cells.Item(15, 25).Value = 25
$ws.Cells.Item(15, 26).Value = 24
$ws.Cells.Item(15, 32).Value = 48
$ws.Cells.Item(15, 38).Value = 38

# Row 17
$ws.Cells.Item(17, 4).Value = 'Renate'
$ws.Cells.Item(17, 5).Value = 'ASD Alcione'
$ws.Cells.Item(17, 6).Value = 2.28
$ws.Cells.Item(17, 7).Value = 2.56
$ws.Cells.Item(17, 8).Value = 3.45
$ws.Cells.Item(17, 9).Value = 4.1
$ws.Cells.Item(17, 10).Value = 2.74
$ws.Cells.Item(17, 11).Value = 3.4
$ws.Cells.Item(17, 14).Value = 1.57
$ws.Cells.Item(17, 16).Value = 1.56
$ws.Cells.Item(17, 17).Value = 2.42
$ws.Cells.Item(17, 19).Value = 2.42
$ws.Cells.Item(17, 22).Value = 1.32
$ws.Cells.Item(17, 23).Value = 1.64

# Row 18
$ws.Cells.Item(18, 4).Value = 'Giana Erminio'
$ws.Cells.Item(18, 5).Value = 'Union Brescia'
$ws.Cells.Item(18, 6).Value = 4.7
$ws.Cells.Item(18, 7).Value = 6.4
$ws.Cells.Item(18, 8).Value = 1.93
$ws.Cells.Item(18, 9).Value = 2.14
$ws.Cells.Item(18, 10).Value = 3
$ws.Cells.Item(18, 11).Value = 3.5
$ws.Cells.Item(18, 14).Value = 1.51
$ws.Cells.Item(18, 15).Value = 1.01
$ws.Cells.Item(18, 16).Value = 1.51
$ws.Cells.Item(18, 17).Value = 2.56
$ws.Cells.Item(18, 19).Value = 2.56
$ws.Cells.Item(18, 22).Value = 1.89
$ws.Cells.Item(18, 23).Value = 1.18

# Row 20
$ws.Cells.Item(20, 35).Value = 110

# Row 21
$ws.Cells.Item(21, 18).Value = 1.25

# Row 22
$ws.Cells.Item(22, 13).Value = 1.14
$ws.Cells.Item(22, 14).Value = 2.38
$ws.Cells.Item(22, 15).Value = 1.64
$ws.Cells.Item(22, 19).Value = 6.2
$ws.Cells.Item(22, 20).Value = 2.24
$ws.Cells.Item(22, 21).Value = 1.64
$ws.Cells.Item(22, 24).Value = 7.8
$ws.Cells.Item(22, 25).Value = 11.5
$ws.Cells.Item(22, 26).Value = 32
$ws.Cells.Item(22, 27).Value = 140
$ws.Cells.Item(22, 28).Value = 6.8
$ws.Cells.Item(22, 29).Value = 7.6
$ws.Cells.Item(22, 30).Value = 21
$ws.Cells.Item(22, 31).Value = 90
$ws.Cells.Item(22, 32).Value = 13
$ws.Cells.Item(22, 33).Value = 13
$ws.Cells.Item(22, 34).Value = 29
$ws.Cells.Item(22, 35).Value = 140
$ws.Cells.Item(22, 36).Value = 34
$ws.Cells.Item(22, 37).Value = 38
$ws.Cells.Item(22, 38).Value = 80
$ws.Cells.Item(22, 39).Value = 290
$ws.Cells.Item(22, 40).Value = 40

# Row 24
$ws.Cells.Item(24, 7).Value = 2.02
$ws.Cells.Item(24, 14).Value = 1.73
$ws.Cells.Item(24, 16).Value = 1.73
$ws.Cells.Item(24, 17).Value = 1.97
$ws.Cells.Item(24, 23).Value = 1.93

# Row 26
$ws.Cells.Item(26, 9).Value = 2.54
$ws.Cells.Item(26, 40).Value = 110

# Row 27
$ws.Cells.Item(27, 26).Value = 970

# Row 29
$ws.Cells.Item(29, 6).Value = 3.25
$ws.Cells.Item(29, 7).Value = 4.2
$ws.Cells.Item(29, 8).Value = 2.4
$ws.Cells.Item(29, 9).Value = 2.8
$ws.Cells.Item(29, 14).Value = 1.55
$ws.Cells.Item(29, 16).Value = 1.55
$ws.Cells.Item(29, 17).Value = 2.32
$ws.Cells.Item(29, 22).Value = 1.56

# Row 30
$ws.Cells.Item(30, 7).Value = 2.08
$ws.Cells.Item(30, 8).Value = 4
$ws.Cells.Item(30, 9).Value = 4.6
$ws.Cells.Item(30, 13).Value = 1.07
$ws.Cells.Item(30, 14).Value = 3.8
$ws.Cells.Item(30, 15).Value = 1.31
$ws.Cells.Item(30, 17).Value = 1.85
$ws.Cells.Item(30, 18).Value = 1.37
$ws.Cells.Item(30, 19).Value = 3.35
$ws.Cells.Item(30, 20).Value = 1.79
$ws.Cells.Item(30, 21).Value = 2.08
$ws.Cells.Item(30, 23).Value = 1.94
$ws.Cells.Item(30, 25).Value = 18
$ws.Cells.Item(30, 34).Value = 22
$ws.Cells.Item(30, 39).Value = 120

# ---- New rows 31-33 (full data, columns A:AO = 41 columns) ----
# NOTE: Date (col B) and Time (col C) values are prefixed with a leading
# apostrophe so Excel stores them as literal text, matching the rest of the
# sheet, instead of auto-converting them to date/time serial values.

$row31 = @('Argentinian Primera Division', "'2026-02-02", "'22:00:00", 'Argentinos Juniors', 'Belgrano', 1.74, 1.8, 6, 7.4, 3.4, 3.6, 1.01, 1.11, 2.72, 1.51, 1.56, 2.46, 1.2, 5, 2.26, 1.64, 1.16, 2.24, 9.800000000000001, 16.5, 55, 280, 6.6, 8.6, 30, 160, 9.6, 11.5, 30, 170, 19.5, 26, 75, 280, 19, 300)
for ($i = 0; $i -lt $row31.Length; $i++) { $ws.Cells.Item(31, $i + 1).Value = $row31[$i] }

$row32 = @('Argentinian Primera Division', "'2026-02-02", "'22:00:00", 'Union Santa Fe', 'Gimnasia Mendoza', 2.04, 2.08, 4.8, 5.3, 3.1, 3.2, 1.01, 1.12, 2.44, 1.61, 1.48, 2.86, 1.16, 6.2, 2.32, 1.66, 1.23, 1.92, 7.8, 12, 36, 150, 6.2, 7.6, 22, 100, 10.5, 12, 30, 170, 26, 32, 75, 330, 34, 180)
for ($i = 0; $i -lt $row32.Length; $i++) { $ws.Cells.Item(32, $i + 1).Value = $row32[$i] }

$row33 = @('Colombian Primera A', "'2026-02-02", "'22:30:00", 'Deportivo Pereira', 'Junior FC Barranquilla', 3.75, 4.5, 2.1, 2.22, 3.25, 3.45, 1.01, 1.01, 2.86, 1.47, 1.62, 2.34, 1.22, 4.5, 1.72, 1.59, 1.82, 1.28, 14.5, 10, 17, 40, 16.5, 10, 16, 38, 42, 25, 30, 75, 1000, 90, 100, 1000, 1000, 1000)
for ($i = 0; $i -lt $row33.Length; $i++) { $ws.Cells.Item(33, $i + 1).Value = $row33[$i] }

